# Swap the contents of data rows 2 and 3 (title, timestamp, historical
# distance, and source uri/hyperlink) so that the two news-item records are
# re-ordered (e.g. after re-running the time-bucket analysis / json import).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current (pre-swap) values --------------------------------
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$e2 = $ws.Range("E2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$e3 = $ws.Range("E3").Value2

# --- remove existing hyperlinks so they can be rebuilt with swapped URLs
$ws.Range("E2").Hyperlinks().Delete()
$ws.Range("E3").Hyperlinks().Delete()

# --- write swapped values back -----------------------------------------
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3
$ws.Range("E2").Value = $e3

$ws.Range("A3").Value = $a2
$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("E3").Value = $e2

# --- re-create hyperlinks on E2/E3 with swapped targets -----------------
$ws.Hyperlinks.Add($ws.Range("E2"), $e3, "", "", $e3)
$ws.Hyperlinks.Add($ws.Range("E3"), $e2, "", "", $e2)
